# Expand the "Positions" sheet of the global multi-asset fixture with two
# additional matched-trade demo rows (rows 5 and 6), extending the used
# range from A1:Z4 to A1:Z6.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Positions")

# --- Row 5: TRD-DEMO1 ------------------------------------------------------
$ws.Cells.Item(5, 1).Value  = "TRD-DEMO1"
$ws.Cells.Item(5, 2).Value  = "LEG-DEMO"
$ws.Cells.Item(5, 3).Value  = "ALLOC-DEMO1"
$ws.Cells.Item(5, 4).Value  = "DEMO-PORT"
$ws.Cells.Item(5, 5).Value  = "ISIN-DEMO1"
$ws.Cells.Item(5, 6).Value  = "GLOBAL"
$ws.Cells.Item(5, 7).Value  = "GMK-APAC"
$ws.Cells.Item(5, 8).Value  = "EQUITY"
$ws.Cells.Item(5, 9).Value  = "USD"
# Text-formatted date fields (stored as literal text, not date serials)
$ws.Cells.Item(5, 10).Value = "'25-Aug-2024"
$ws.Cells.Item(5, 11).Value = "'27-Aug-2024"
$ws.Cells.Item(5, 12).Value = 0.0092
$ws.Cells.Item(5, 13).Value = 100000
$ws.Cells.Item(5, 14).Value = 99500
$ws.Cells.Item(5, 15).Value = 11
$ws.Cells.Item(5, 16).Value = 175
$ws.Cells.Item(5, 17).Value = 1500
$ws.Cells.Item(5, 18).Value = 166.4
$ws.Cells.Item(5, 19).Value = "BOOKED"
$ws.Cells.Item(5, 20).Value = "AFFIRMED"
$ws.Cells.Item(5, 21).Value = "DemoCounterparty"
$ws.Cells.Item(5, 22).Value = "DemoBroker"
$ws.Cells.Item(5, 23).Value = "Y"
$ws.Cells.Item(5, 24).Value = "GLOBAL-DEMO"
# Confidence seed stored as text (not numeric) in the source fixture
$ws.Cells.Item(5, 25).Value = "'0.99"
$ws.Cells.Item(5, 26).Value = "MATCH-CONFIRMED"

# --- Row 6: TRD-DEMO2 ------------------------------------------------------
$ws.Cells.Item(6, 1).Value  = "TRD-DEMO2"
$ws.Cells.Item(6, 2).Value  = "LEG-DEMO"
$ws.Cells.Item(6, 3).Value  = "ALLOC-DEMO2"
$ws.Cells.Item(6, 4).Value  = "DEMO-PORT"
$ws.Cells.Item(6, 5).Value  = "ISIN-DEMO2"
$ws.Cells.Item(6, 6).Value  = "APAC"
$ws.Cells.Item(6, 7).Value  = "GMK-APAC"
$ws.Cells.Item(6, 8).Value  = "EQUITY"
$ws.Cells.Item(6, 9).Value  = "JPY"
$ws.Cells.Item(6, 10).Value = "'21-Aug-2024"
$ws.Cells.Item(6, 11).Value = "'23-Aug-2024"
$ws.Cells.Item(6, 12).Value = 0.0092
$ws.Cells.Item(6, 13).Value = 260000
$ws.Cells.Item(6, 14).Value = 259600
$ws.Cells.Item(6, 15).Value = 11
$ws.Cells.Item(6, 16).Value = 180
$ws.Cells.Item(6, 17).Value = 1500
$ws.Cells.Item(6, 18).Value = 173.3
$ws.Cells.Item(6, 19).Value = "BOOKED"
$ws.Cells.Item(6, 20).Value = "AFFIRMED"
$ws.Cells.Item(6, 21).Value = "PrimeAPAC"
$ws.Cells.Item(6, 22).Value = "BrokerAPAC"
$ws.Cells.Item(6, 23).Value = "Y"
$ws.Cells.Item(6, 24).Value = "GLOBAL-DEMO"
$ws.Cells.Item(6, 25).Value = "'0.95"
$ws.Cells.Item(6, 26).Value = "Positions"
